$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("TblNames")

# Insert a new row at row 19 (pushes Tbl_DailyExposure.. down by one)
$ws.Rows.Item(19).Insert()

# New row content: Tbl_CustomersList with ID 727
$ws.Cells.Item(19, 1).Value = 727
$ws.Cells.Item(19, 2).Value = "Tbl_CustomersList"

# Renumber every data row's ID (column A) to the new sequence 710..760
$names = @(
  "CreditCheckFailureTXTFileAmended",
  "Indirizzi",
  "Paste Errors",
  "Tbl_AdditionalQueryData",
  "Tbl_Areas",
  "Tbl_Banks",
  "Tbl_Cash_Target",
  "Tbl_Cash_Target_Breakdown",
  "Tbl_CashCollected",
  "Tbl_Channels",
  "Tbl_CL",
  "Tbl_Countries",
  "Tbl_Country_Internal_Contact",
  "Tbl_credit_check_failures",
  "Tbl_Currencies",
  "Tbl_Customer_Status",
  "Tbl_Customers",
  "Tbl_CustomersList",
  "Tbl_DailyExposure",
  "Tbl_Deductions",
  "Tbl_DepartmentNames",
  "Tbl_DocumentsToBeErased",
  "Tbl_EmailAddresses",
  "Tbl_EmailSoftware",
  "Tbl_GeneralChart",
  "Tbl_HelpPages",
  "Tbl_Historical_Statements",
  "Tbl_HistoricalCLsAndStatements",
  "Tbl_HoldTypesToBeConsideredForCreditCheckReleases",
  "Tbl_InvoiceAttachments",
  "Tbl_Invoices",
  "Tbl_Invoices_History",
  "Tbl_Languages",
  "Tbl_Link_Customer_Internal_Email_Address",
  "TBL_LinkTemplateEmailAddress",
  "Tbl_LOGICollectionsManagementReport",
  "Tbl_MonthEnd",
  "Tbl_NEWS",
  "Tbl_PaymentData",
  "Tbl_queries",
  "Tbl_RSS",
  "Tbl_SmartActivities",
  "Tbl_Templates",
  "Tbl_Timezones",
  "Tbl_Top5ComplainingCustomers",
  "Tbl_Top5WorseCustomers",
  "Tbl_Types",
  "Tbl_Users",
  "Tbl_WhoPaidYesterday",
  "TblGeneral",
  "TblMain"
)

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $i + 2
  $id = 710 + $i
  $ws.Cells.Item($row, 1).Value = $id
  $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Update the sheet's used-range dimension and the workbook-level defined name
$wb.Names.Item("TblNames").RefersTo = "='TblNames'!`$A`$1:`$B`$52"

Write-Output "ok"
